# Update the "Metadata" sheet of the CodeSystem-wh-payer-indicators workbook:
# bump the IG version/date, set the real publisher, replace the stray
# "Contact / No display for ContactDetail" row with a single "Jurisdiction"
# row, and record "Case Sensitive" = true.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" row (old row 11); this shifts rows 12-22 up to 11-21.
$ws.Rows.Item(11).Delete()

# Update Version value (row 3)
$ws.Range("B3").Value = "6.0.0"

# Update Date value (row 8)
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Update Publisher value (row 9)
$ws.Range("B9").Value = "Alvearie Team"

# Replace old duplicate "Contact" row (row 10) with the new "Jurisdiction" row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Update Case Sensitive value (row 14 after the row shift) - must be stored
# as literal text "true", not an auto-coerced Boolean, so enter it as a
# formula returning the text and flatten it back to a literal via paste-values.
$caseCell = $ws.Range("B14")
$caseCell.Formula = "=""true"""
$caseCell.Copy()
$caseCell.PasteSpecial(-4163)
